# Update workflow control table with new BSU, NVIS and Marine input paths
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2 (BSU): replace placeholder "NA" with the new BSU output path and
#     turn it into a hyperlink, vertically centered like the other link cells ---
$ws.Range("B2").VerticalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("B2"), "file:///\\fs1-cbr.nexus.csiro.au\%7bev-neap%7d\work\BSU\outputs\BSU_NEAP\BSU_NEAP_epsg3577_250m.tif") | Out-Null
$ws.Range("B2").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\BSU\outputs\BSU_NEAP\BSU_NEAP_epsg3577_250m.tif"

# --- Row 3 (Marine): point RawDataPath at the new NESP-MERI benthic raster ---
$ws.Range("B3").Font.Name = "Aptos"
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("B3").Value = "//fs1-cbr.nexus.csiro.au/{ev-neap}/work/extent/inputs/raw/Marine/NVE-Benthic/NESP-MERI_Natural_Values_Ecosystems_withVAT.tif"

# --- Row 4 (NVIS_NEAP): point RawDataPath at the new pre-1750 NVIS raster ---
$ws.Range("B4").Font.Name = "Aptos"
$ws.Range("B4").Value = "//fs1-cbr.nexus.csiro.au/{ev-neap}/work/extent/processing/NEAP_intermediate/NVIS_PRE1750_IUCNGET_DK_20240714.tif"

# --- Update the visible selection to span the DatasetName/RawDataPath rows ---
$ws.Range("B1:B9").Select()
